$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

Set-TextValue "D2" "30.806.27"
$ws.Range("E2").Value = "  +2.34%  "

Set-TextValue "D3" "2.112.72"
$ws.Range("E3").Value = "  +9.04%  "

Set-TextValue "D4" "1.002"
$ws.Range("E4").Value = "  +0.20%  "

Set-TextValue "D5" "333.70"
$ws.Range("E5").Value = "  +4.22%  "

$ws.Range("E6").Value = "  +0.25%  "

Set-TextValue "D7" "0.5269"
$ws.Range("E7").Value = "  +3.66%  "

Set-TextValue "D8" "0.4381"
$ws.Range("E8").Value = "  +8.36%  "

Set-TextValue "D9" "0.09022"
$ws.Range("E9").Value = "  +7.28%  "

Set-TextValue "D10" "45.96"
$ws.Range("E10").Value = "  +8.74%  "

Set-TextValue "D11" "1.177"
$ws.Range("E11").Value = "  +4.90%  "

$ws.Range("E12").Value = "  +3.79%  "

Set-TextValue "D13" "2.110.85"
$ws.Range("E13").Value = "  +9.58%  "

Set-TextValue "D14" "6.766"
$ws.Range("E14").Value = "  +5.34%  "

Set-TextValue "D15" "7.817"
$ws.Range("E15").Value = "  +7.17%  "

Set-TextValue "D16" "97.29"
$ws.Range("E16").Value = "  +4.60%  "

Set-TextValue "D17" "1.004"
$ws.Range("E17").Value = "  +0.09%  "

Set-TextValue "D18" "0.00001129"
$ws.Range("E18").Value = "  +2.75%  "

Set-TextValue "D19" "0.06675"
$ws.Range("E19").Value = "  +2.55%  "

$ws.Range("E20").Value = "  +2.85%  "

$ws.Range("E21").Value = "  +0.15%  "

Set-TextValue "D22" "6.354"
$ws.Range("E22").Value = "  +6.34%  "

Set-TextValue "D23" "30.862.09"
$ws.Range("E23").Value = "  +2.50%  "

$ws.Range("E24").Value = "  +6.91%  "

Set-TextValue "D25" "2.359.58"
$ws.Range("E25").Value = "  +10.07%  "

Set-TextValue "D26" "2.266"
$ws.Range("E26").Value = "  +3.51%  "

Set-TextValue "D27" "22.78"
$ws.Range("E27").Value = "  +3.28%  "

Set-TextValue "D28" "2.560"
$ws.Range("E28").Value = "  +12.15%  "

Set-TextValue "D29" "162.52"
$ws.Range("E29").Value = "  -0.16%  "

Set-TextValue "D30" "132.98"
$ws.Range("E30").Value = "  +2.45%  "

Set-TextValue "D31" "1.169"
$ws.Range("E31").Value = "  +3.22%  "

Set-TextValue "D32" "0.1071"
$ws.Range("E32").Value = "  +2.30%  "

Set-TextValue "D33" "6.229"
$ws.Range("E33").Value = "  +3.98%  "

Set-TextValue "D35" "1.542"
$ws.Range("E35").Value = "  +22.24%  "

Set-TextValue "D36" "0.02602"
$ws.Range("E36").Value = "  +5.93%  "

Set-TextValue "D37" "5.540"
$ws.Range("E37").Value = "  +4.03%  "

Set-TextValue "D38" "0.06736"
$ws.Range("E38").Value = "  +4.11%  "

Set-TextValue "D39" "9.510"
$ws.Range("E39").Value = "  +9.08%  "

Set-TextValue "D40" "12.73"
$ws.Range("E40").Value = "  +8.57%  "

Set-TextValue "D41" "0.2270"
$ws.Range("E41").Value = "  +5.32%  "

Set-TextValue "D42" "0.6835"
$ws.Range("E42").Value = "  +5.19%  "

Set-TextValue "D43" "1.249"
$ws.Range("E43").Value = "  +2.23%  "

Set-TextValue "D44" "0.6460"
$ws.Range("E44").Value = "  +6.29%  "

Set-TextValue "D45" "1.002"
$ws.Range("E45").Value = "  +0.24%  "

Set-TextValue "D46" "14.04"
$ws.Range("E46").Value = "  +5.36%  "

Set-TextValue "D47" "2.228"
$ws.Range("E47").Value = "  +2.25%  "

Set-TextValue "D48" "3.672"
$ws.Range("E48").Value = "  +1.23%  "

Set-TextValue "D49" "1.275"
$ws.Range("E49").Value = "  +5.19%  "

Set-TextValue "D50" "82.39"
$ws.Range("E50").Value = "  +5.34%  "

Set-TextValue "D51" "0.07108"
$ws.Range("E51").Value = "  +3.98%  "
